$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.230.91'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").Value = '2.244.05'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.45'
$ws.Range("E5").Value = '  -2.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.08'
$ws.Range("E6").Value = '  -3.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.572'
$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("E8").Value = '  +0.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.522'
$ws.Range("E9").Value = '  -1.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.69'
$ws.Range("E10").Value = '  -3.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0809'
$ws.Range("E11").Value = '  -1.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.19'
$ws.Range("E12").Value = '  -2.10%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '2.338.19'
$ws.Range("E14").Value = '  +4.70%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.585.86'
$ws.Range("E15").Value = '  +0.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.830'
$ws.Range("E16").Value = '  -0.93%  '

$ws.Range("E17").Value = '  -2.64%  '

$ws.Range("D18").Value = '43.984.24'
$ws.Range("E18").Value = '  +0.53%  '

$ws.Range("E19").Value = '  -0.77%  '

$ws.Range("E20").Value = '  +1.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.08'
$ws.Range("E21").Value = '  -7.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.55'
$ws.Range("E22").Value = '  -0.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.58'
$ws.Range("E23").Value = '  +0.53%  '

$ws.Range("E24").Value = '  -1.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.99'
$ws.Range("E25").Value = '  -1.69%  '

$ws.Range("E26").Value = '  +0.10%  '

$ws.Range("E27").Value = '  -1.70%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.20'
$ws.Range("E28").Value = '  +3.10%  '

$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.89'
$ws.Range("E29").Value = '  +3.74%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.03'
$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("E31").Value = '  -2.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.13'
$ws.Range("E32").Value = '  -2.08%  '

$ws.Range("E33").Value = '  -4.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.61'
$ws.Range("E34").Value = '  -1.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.21'
$ws.Range("E35").Value = '  -3.42%  '

$ws.Range("E36").Value = '  +2.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.107'
$ws.Range("E37").Value = '  -1.41%  '

$ws.Range("E38").Value = '  -7.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.54'
$ws.Range("E39").Value = '  +0.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.82'
$ws.Range("E40").Value = '  -4.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.19'
$ws.Range("E41").Value = '  -9.02%  '

$ws.Range("E42").Value = '  -2.95%  '

$ws.Range("E43").Value = '  +0.30%  '

$ws.Range("D44").Value = '1.747.69'
$ws.Range("E44").Value = '  +2.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '82.32'
$ws.Range("E45").Value = '  -0.50%  '

$ws.Range("E46").Value = '  -2.19%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.10'
$ws.Range("E47").Value = '  +0.52%  '

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '99.76'
$ws.Range("E48").Value = '  -1.87%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.93'
$ws.Range("E49").Value = '  -4.28%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.09'
$ws.Range("E50").Value = '  -0.62%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.57'
$ws.Range("E51").Value = '  -2.73%  '
